$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gained two more year columns: N = 2022, O = 2023.
# Copy the formatting from the previous year's column (M) so the new
# cells pick up the same style index used throughout the header/data rows.
$ws.Range("M4").Copy()
$ws.Range("N4:O4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N4").Value = 2022
$ws.Range("O4").Value = 2023

$ws.Range("M5").Copy()
$ws.Range("N5:O5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N5").Value = 6.53
$ws.Range("O5").Value = 6.53

$excel.CutCopyMode = $false

# A handful of rows got their height normalised to match the new layout.
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 17.25

# Reset the saved selection to A1 instead of the stray P6 left in the
# source file.
$ws.Range("A1").Select()
